$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at the very bottom of the data block first
# (row 22) rather than at row 2 directly -- inserting right under the bold
# header row would copy its formatting down onto the new row, polluting
# styles.xml with an extra cellXf that the target workbook doesn't have.
$ws.Rows("22:22").Insert()

# Shift the existing data rows 2..21 down into 3..22, bottom row first, by
# copying whole rows (this preserves each cell's text/number typing, so
# shared-string cells remain shared strings instead of being re-typed).
for ($r = 21; $r -ge 2; $r--) {
    $srcRow = $ws.Range("A" + $r + ":Y" + $r)
    $dstRow = $ws.Range("A" + ($r + 1) + ":Y" + ($r + 1))
    $srcRow.Copy($dstRow)
}

# Row 2 is now free for the new "iBeam Technology" (아이빔테크놀로지) IPO entry
# pushed in this dataset update. A handful of the text columns (the two
# demand-forecast dates, the listing date, and the subscription-competition
# percentage) look like a date/number to Excel's type inference, so they
# need to be pre-marked as Text before assignment or they'd silently turn
# into date serials / percentages instead of staying literal strings.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("O2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-07-15"
$ws.Range("B2").Value = "2024-07-19"
$ws.Range("C2").Value = "2024-08-06"
$ws.Range("D2").Value = "삼성"
$ws.Range("E2").Value = "아이빔테크놀로지"
$ws.Range("F2").Value = 2234000
$ws.Range("G2").Value = 2234000
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 7300
$ws.Range("J2").Value = 8500
$ws.Range("K2").Value = 14965620
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 10000
$ws.Range("N2").Value = "1011.50:1"
$ws.Range("O2").Value = "0.49%"
$ws.Range("P2").Value = 1337894447
$ws.Range("Q2").Value = 4495088983
$ws.Range("R2").Value = 99398351
$ws.Range("S2").Value = -3330678562
$ws.Range("T2").Value = -2905069801
$ws.Range("U2").Value = -1314068310
$ws.Range("V2").Value = -9735376868
$ws.Range("W2").Value = -5043515600
$ws.Range("X2").Value = -1257080998
$ws.Range("Y2").Value = "생체현미경, CRO 서비스"
